# Clean up form formatting
#
# - "survey" sheet field A8 was mis-typed as the ODK type "string" (not a
#   valid type in the dropdown's validation list); correct it to "text".
# - The "settings" sheet was left active/selected; switch focus back to the
#   "survey" sheet, which becomes the active tab, with cell D17 selected in
#   the frozen bottom-right pane (rather than the previous D26).

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# Fix the mistyped field type in the survey sheet.
$survey.Range("A8").Value = "text"

# Restore focus: settings was the active/selected tab before, move back to
# survey and leave D17 as the selected cell (matches the author's final
# on-screen selection).
$settings.Activate() | Out-Null
$survey.Activate() | Out-Null
$survey.Range("D17").Select() | Out-Null
